# Banco de dados alimentados com registros
# Populates the phrasal-verb bank (Sheet 1) with 9 new entries, and refreshes
# sheet-level view/print/column settings to match the saved workbook state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column E width: widen to a fixed 31 characters (drops the old auto best-fit) ---
$ws.Columns.Item(5).ColumnWidth = 30.1

# --- Print setup: A4 paper, portrait orientation ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- New phrasal-verb rows (id_phrasal, titulo, traducaoTitulo, descricao, exemplo, exemploTraducao, dificuldade) ---

# Row 3: Act up
$ws.Cells.Item(3,1).NumberFormat = "@"
$ws.Cells.Item(3,1).Value = "2"
$ws.Cells.Item(3,1).Style = "Normal"
$ws.Cells.Item(3,2).Value = "Act up"
$ws.Cells.Item(3,3).Value = "Agir"
$ws.Cells.Item(3,4).Value = "Comportar-se mal, não funcionar como deveria"
$ws.Cells.Item(3,5).Value = "The babysitter had a difficult time, the children acted up all evening."
$ws.Cells.Item(3,6).Value = "A babá teve um tempo difícil, as crianças agiram a noite toda."
$ws.Cells.Item(3,7).Value = 2

# Row 4: Ask around
$ws.Cells.Item(4,1).NumberFormat = "@"
$ws.Cells.Item(4,1).Value = "3"
$ws.Cells.Item(4,1).Style = "Normal"
$ws.Cells.Item(4,2).Value = "Ask around"
$ws.Cells.Item(4,3).Value = "Pergunte por aí"
$ws.Cells.Item(4,4).Value = "Convidar alguém para sair"
$ws.Cells.Item(4,5).Value = "He asked her out last night."
$ws.Cells.Item(4,6).Value = "Ele a chamou para sair ontem a noite"
$ws.Cells.Item(4,7).Value = 2

# Row 5: Blow up
$ws.Cells.Item(5,1).NumberFormat = "@"
$ws.Cells.Item(5,1).Value = "4"
$ws.Cells.Item(5,1).Style = "Normal"
$ws.Cells.Item(5,2).Value = "Blow up"
$ws.Cells.Item(5,3).Value = "Explodir"
$ws.Cells.Item(5,4).Value = "Explodir (bomba, etc.); estourar (tempestade, escândalo, crise, guerra, etc.); perder a paciência"
$ws.Cells.Item(5,5).Value = "There was a huge bang as the fuel tank blew up."
$ws.Cells.Item(5,6).Value = "Houve uma explosão enorme quando o tanque de combustível explodiu."
$ws.Cells.Item(5,7).Value = 2

# Row 6: Break down
$ws.Cells.Item(6,1).NumberFormat = "@"
$ws.Cells.Item(6,1).Value = "5"
$ws.Cells.Item(6,1).Style = "Normal"
$ws.Cells.Item(6,2).Value = "Break down"
$ws.Cells.Item(6,3).Value = "demolir"
$ws.Cells.Item(6,4).Value = "Quebrar, parar de funcionar"
$ws.Cells.Item(6,5).Value = "Our car broke down on the road."
$ws.Cells.Item(6,6).Value = "Nosso carro quebrou na estrada."
$ws.Cells.Item(6,7).Value = 2

# Row 7: Break up
$ws.Cells.Item(7,1).NumberFormat = "@"
$ws.Cells.Item(7,1).Value = "6"
$ws.Cells.Item(7,1).Style = "Normal"
$ws.Cells.Item(7,2).Value = "Break up"
$ws.Cells.Item(7,3).Value = "Rompimento"
$ws.Cells.Item(7,4).Value = "Separar-se, terminar um relacionamento, uma parceria, etc."
$ws.Cells.Item(7,5).Value = "They broke up four years ago."
$ws.Cells.Item(7,6).Value = "Eles terminaram há quatro anos."
$ws.Cells.Item(7,7).Value = 2

# Row 8: Cheer up
$ws.Cells.Item(8,1).NumberFormat = "@"
$ws.Cells.Item(8,1).Value = "7"
$ws.Cells.Item(8,1).Style = "Normal"
$ws.Cells.Item(8,2).Value = "Cheer up"
$ws.Cells.Item(8,3).Value = "Anime-se"
$ws.Cells.Item(8,4).Value = "Animar-se"
$ws.Cells.Item(8,5).Value = "She cheered up when he got home."
$ws.Cells.Item(8,6).Value = "Ela se animou quando chegou em casa."
$ws.Cells.Item(8,7).Value = 3

# Row 9: Chip in
$ws.Cells.Item(9,1).NumberFormat = "@"
$ws.Cells.Item(9,1).Value = "8"
$ws.Cells.Item(9,1).Style = "Normal"
$ws.Cells.Item(9,2).Value = "Chip in"
$ws.Cells.Item(9,3).Value = "Sem Tradução"
$ws.Cells.Item(9,4).Value = "Contribuir com dinheiro, ""fazer uma vaquinha"""
$ws.Cells.Item(9,5).Value = "If everyone chips in we'll be able to buy her a nice present."
$ws.Cells.Item(9,6).Value = "Se todos mexerem, poderemos comprá-lo um bom presente."
$ws.Cells.Item(9,7).Value = 3

# Row 10: Come apart
$ws.Cells.Item(10,1).NumberFormat = "@"
$ws.Cells.Item(10,1).Value = "9"
$ws.Cells.Item(10,1).Style = "Normal"
$ws.Cells.Item(10,2).Value = "Come apart"
$ws.Cells.Item(10,3).Value = "Separar"
$ws.Cells.Item(10,4).Value = "Separar-se, desfazer-se em pedaços, quebrar, desmoronar"
$ws.Cells.Item(10,5).Value = "The top and the bottom come apart if you pull hard enough."
$ws.Cells.Item(10,6).Value = "O topo e o fundo se aparecem se você puxar o suficiente."
$ws.Cells.Item(10,7).Value = 2

# Row 11: Come off
$ws.Cells.Item(11,1).NumberFormat = "@"
$ws.Cells.Item(11,1).Value = "10"
$ws.Cells.Item(11,1).Style = "Normal"
$ws.Cells.Item(11,2).Value = "Come off"
$ws.Cells.Item(11,3).Value = "Saia"
$ws.Cells.Item(11,4).Value = "Sair, desaparecer"
$ws.Cells.Item(11,5).Value = "That mark on your dress won't come off."
$ws.Cells.Item(11,6).Value = "Essa marca no seu vestido não vai sair."
$ws.Cells.Item(11,7).Value = 3

# --- Select the whole sheet (Ctrl+A), mirroring the selection saved in the file ---
$ws.Cells.Select() | Out-Null
